$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Air India'
$ws.Range("B2").Value = '09:55'
$ws.Range("C2").Value = '03 h 05 m'
$ws.Range("D2").Value = '₹ 8,248'

$ws.Range("A3").Value = 'Air India'
$ws.Range("B3").Value = '20:15'
$ws.Range("C3").Value = '02 h 30 m'
$ws.Range("D3").Value = '₹ 8,248'

$ws.Range("A4").Value = 'SpiceJet'
$ws.Range("B4").Value = '21:40'
$ws.Range("C4").Value = '02 h 45 m'
$ws.Range("D4").Value = '₹ 8,982'

$ws.Range("A5").Value = 'SpiceJet'
$ws.Range("B5").Value = '08:30'
$ws.Range("C5").Value = '05 h 30 m'
$ws.Range("D5").Value = '₹ 8,982'

$ws.Range("A6").Value = 'Air India'
$ws.Range("B6").Value = '06:00'
$ws.Range("C6").Value = '02 h 55 m'
$ws.Range("D6").Value = '₹ 8,983'

$ws.Range("A7").Value = 'IndiGo'
$ws.Range("B7").Value = '06:20'
$ws.Range("C7").Value = '02 h 50 m'
$ws.Range("D7").Value = '₹ 8,983'

$ws.Range("A8").Value = 'Vistara'
$ws.Range("B8").Value = '07:05'
$ws.Range("C8").Value = '02 h 50 m'
$ws.Range("D8").Value = '₹ 8,983'

$ws.Range("A9").Value = 'IndiGo'
$ws.Range("B9").Value = '08:45'
$ws.Range("C9").Value = '02 h 45 m'
$ws.Range("D9").Value = '₹ 8,983'

$ws.Range("A10").Value = 'Vistara'
$ws.Range("B10").Value = '10:35'
$ws.Range("C10").Value = '02 h 50 m'
$ws.Range("D10").Value = '₹ 8,983'

$ws.Range("A11").Value = 'IndiGo'
$ws.Range("B11").Value = '10:40'
$ws.Range("C11").Value = '02 h 45 m'
$ws.Range("D11").Value = '₹ 8,983'

$ws.Range("A12").Value = 'IndiGo'
$ws.Range("B12").Value = '13:20'
$ws.Range("C12").Value = '03 h'
$ws.Range("D12").Value = '₹ 8,983'

$ws.Range("A13").Value = 'IndiGo'
$ws.Range("B13").Value = '15:10'
$ws.Range("C13").Value = '02 h 50 m'
$ws.Range("D13").Value = '₹ 8,983'

$ws.Range("A14").Value = 'IndiGo'
$ws.Range("B14").Value = '16:35'
$ws.Range("C14").Value = '02 h 50 m'
$ws.Range("D14").Value = '₹ 8,983'

$ws.Range("A15").Value = 'Air India'
$ws.Range("B15").Value = '16:55'
$ws.Range("C15").Value = '02 h 55 m'
$ws.Range("D15").Value = '₹ 8,983'

$ws.Range("A16").Value = 'Vistara'
$ws.Range("B16").Value = '17:15'
$ws.Range("C16").Value = '02 h 45 m'
$ws.Range("D16").Value = '₹ 8,983'

$ws.Range("A17").Value = 'IndiGo'
$ws.Range("B17").Value = '18:15'
$ws.Range("C17").Value = '02 h 50 m'
$ws.Range("D17").Value = '₹ 8,983'

$ws.Range("A18").Value = 'IndiGo'
$ws.Range("B18").Value = '19:35'
$ws.Range("C18").Value = '02 h 50 m'
$ws.Range("D18").Value = '₹ 8,983'

$ws.Range("A19").Value = 'IndiGo'
$ws.Range("B19").Value = '22:50'
$ws.Range("C19").Value = '02 h 55 m'
$ws.Range("D19").Value = '₹ 8,983'

$ws.Range("A20").Value = 'IndiGo'
$ws.Range("B20").Value = '11:00'
$ws.Range("C20").Value = '04 h 35 m'
$ws.Range("D20").Value = '₹ 9,036'

$ws.Range("A21").Value = 'SpiceJet'
$ws.Range("B21").Value = '19:40'
$ws.Range("C21").Value = '13 h 10 m'
$ws.Range("D21").Value = '₹ 13,981'

$ws.Range("A22").Value = 'SpiceJet'
$ws.Range("B22").Value = '19:40'
$ws.Range("C22").Value = '13 h 10 m'
$ws.Range("D22").Value = '₹ 13,509'

